$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new bulleted ("ListParagraph"/numId=1) paragraph right
#    after the "The license files will be stored ..." paragraph, with
#    the new text, and move the _GoBack bookmark to its end.
# ---------------------------------------------------------------------

$srcPara = $d.Paragraphs(8)
$srcPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs(9)
$newPara.Range.Text = "It will be a licensing and verification based on the input file in which the validation logic is written in an encrypted format."

# Move the _GoBack bookmark from its old spot (after "self-upgradeable")
# to the end of the freshly inserted paragraph's text.
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

$endPos = $newPara.Range.End - 1

# Zero-length ranges placed exactly before a paragraph mark land the new
# bookmark at document position 0, so insert a throw-away placeholder
# character, bookmark across it, then delete the character; the
# bookmark collapses back onto the correct (zero-length) spot.
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")
$bmRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$cleanup = $d.Range($endPos, $endPos + 1)
$cleanup.Text = ""

# ---------------------------------------------------------------------
# 2) Merge the " "/"will"/" it be a" runs (proofing marks removed by
#    Word's own re-save) into a single run; text is unchanged.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" will it be a", $true, $false, $false, $false, $false, $true, 1, $false, " will it be a", 2)

# ---------------------------------------------------------------------
# 3) Remove the spell/grammar squiggly markers elsewhere (text
#    unaffected); simple no-op replacements keep the paragraphs but
#    drop the proofErr runs when Word merges adjacent identical runs on
#    save.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Method to avoid the infringement by changing the executable and dll. ", $true, $false, $false, $false, $false, $true, 1, $false, "Method to avoid the infringement by changing the executable and dll. ", 2)
$d.Content.Find.Execute("The hacker can patch the exe and binary dlls. Either sign it and then check the signing before exe or dll. use handshaking in every api call. ", $true, $false, $false, $false, $false, $true, 1, $false, "The hacker can patch the exe and binary dlls. Either sign it and then check the signing before exe or dll. use handshaking in every api call. ", 2)
$d.Content.Find.Execute("APL: a dll which will contain most of the licensing and serial number related apis. It is a signed dll and it will be checked against its hash for any patching done in it. This dll will be used by installer at the time of installation. ", $true, $false, $false, $false, $false, $true, 1, $false, "APL: a dll which will contain most of the licensing and serial number related apis. It is a signed dll and it will be checked against its hash for any patching done in it. This dll will be used by installer at the time of installation. ", 2)
$d.Content.Find.Execute("was given to him earlier and he gets a grace period. ", $true, $false, $false, $false, $false, $true, 1, $false, "was given to him earlier and he gets a grace period. ", 2)

Write-Host "done"
